# Add a new "2022-Q3" quarterly sheet (as the second sheet, right after "总计"
# and before "2022-Q2"), populate it with fund holding data, and add a
# corresponding summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

$total = $wb.Worksheets.Item(1)      # "总计" sheet
$q2Sheet = $wb.Worksheets.Item(2)    # current "2022-Q2" sheet (will stay 2nd data sheet)

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet right before the "2022-Q2" sheet
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# A pristine, never-written cell on the new sheet - used purely as a
# "format donor" so we can strip any unwanted style (e.g. leftover from
# setting NumberFormat) off other cells via PasteSpecial(formats).
$clean = $newSheet.Cells.Item(100,100)

# Copy header-row formatting (bold / bordered / centered style) from the
# "总计" sheet header onto the new sheet's header cells.
$total.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

# Copy the numeric "index" column formatting (style used for column A data)
$total.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial($xlPasteFormats)

# Header row text
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Source data: columns B,D,E,F,G are stored as plain TEXT in the workbook
# (even though several look numeric, e.g. fund codes with leading zeros,
# or percentages), column H is a real number, and column A is a
# zero-based numeric row index. Column C (fund name) is plain text too,
# but has no digits-only content so it is safe to assign directly.
$data = @(
    @(0, "011868", "中信建投远见回报混合A",   "6.14", "94.99", "2.49", "0.1529", 9),
    @(1, "011869", "中信建投远见回报混合C",   "1.39", "94.99", "2.49", "0.0346", 9),
    @(2, "166109", "信澳量化先锋混合（LOF）A", "0.79", "88.99", "2.48", "0.0196", 7),
    @(3, "000398", "华富灵活配置混合",         "0.12", "94.04", "3.01", "0.0036", 8),
    @(4, "166110", "信澳量化先锋混合（LOF）C", "0.11", "88.99", "2.48", "0.0027", 7)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r,1).Value = $row[0]

    # Text columns: force text format *before* assigning so numeric-looking
    # strings (fund codes, percentages) keep leading zeros / exact text.
    $newSheet.Cells.Item($r,2).NumberFormat = "@"
    $newSheet.Cells.Item($r,2).Value = $row[1]

    $newSheet.Cells.Item($r,3).Value = $row[2]

    $newSheet.Cells.Item($r,4).NumberFormat = "@"
    $newSheet.Cells.Item($r,4).Value = $row[3]

    $newSheet.Cells.Item($r,5).NumberFormat = "@"
    $newSheet.Cells.Item($r,5).Value = $row[4]

    $newSheet.Cells.Item($r,6).NumberFormat = "@"
    $newSheet.Cells.Item($r,6).Value = $row[5]

    $newSheet.Cells.Item($r,7).NumberFormat = "@"
    $newSheet.Cells.Item($r,7).Value = $row[6]

    $newSheet.Cells.Item($r,8).Value = $row[7]

    # Strip the style index that leaks in from setting NumberFormat="@" on
    # B,D,E,F,G (those cells should stay visually unstyled, just like the
    # other quarter sheets).
    $clean.Copy()
    $newSheet.Range("B" + $r + ":G" + $r).PasteSpecial($xlPasteFormats)

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Insert a new row 2 on the "总计" sheet for the 2022-Q3 summary
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# Re-apply the original column-A style (it is lost on row insert) by
# copying formats from the row just below (still carrying the old style).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial($xlPasteFormats)

# Columns B-D are unstyled (default) data cells; clear the leftover
# "inserted row" formatting the same way, using a row that still has the
# original (unstyled) formatting as the source.
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial($xlPasteFormats)

# Rewrite the whole table (5 quarters) so column A keeps its 0..4 running
# index while columns B-D hold the correct quarter/count/value for each
# row (the newest quarter, 2022-Q3, goes on top).
$totalData = @(
    @(0, "2022-Q3", 5, 0.21),
    @(1, "2022-Q2", 4, 0.07000000000000001),
    @(2, "2022-Q1", 9, 0.67),
    @(3, "2021-Q4", 3, 0.13),
    @(4, "2021-Q1", 1, 0.09)
)

$r = 2
foreach ($row in $totalData) {
    $total.Cells.Item($r,1).Value = $row[0]
    $total.Cells.Item($r,2).Value = $row[1]
    $total.Cells.Item($r,3).Value = $row[2]
    $total.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}
